$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 555.4375
$ws.Range("I2").Value = 347.25
$ws.Range("K2").Value = 347.25
$ws.Range("M2").Value = -234.25
$ws.Range("H100").Value = 697.7143
$ws.Range("I100").Value = 496.25
$ws.Range("J100").Value = 966.3333
$ws.Range("K100").Value = 496.25
$ws.Range("L100").Value = 966.3333
$ws.Range("M100").Value = 44.75
$ws.Range("N100").Value = -2048.3333
$ws.Range("H101").Value = 45455696
$ws.Range("J101").Value = 3660
$ws.Range("L101").Value = 10980
$ws.Range("N101").Value = -14224
$ws.Range("H116").Value = 4860.4375
$ws.Range("I116").Value = 3899
$ws.Range("J116").Value = 5608.222
$ws.Range("K116").Value = 3899
$ws.Range("L116").Value = 5608.222
$ws.Range("M116").Value = -457
$ws.Range("N116").Value = -12492.222
$ws.Range("H132").Value = 20836066
$ws.Range("I132").Value = 21279374
$ws.Range("J132").Value = 600
$ws.Range("K132").Value = 63838122
$ws.Range("L132").Value = 1800
$ws.Range("M132").Value = -63835592
$ws.Range("N132").Value = -6860
$ws.Range("H138").Value = 3546.7021
$ws.Range("I138").Value = 3092.4285
$ws.Range("K138").Value = 9277.2855
$ws.Range("M138").Value = -4137.2855

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 12194.706
$ws.Range("I32").Value = 8571.929
$ws.Range("J32").Value = 16605.043
$ws.Range("K32").Value = 8571.929
$ws.Range("L32").Value = 16605.043
$ws.Range("M32").Value = -8284.929
$ws.Range("N32").Value = -17179.043
$ws.Range("H45").Value = 4798436.5
$ws.Range("J45").Value = 7558.125
$ws.Range("L45").Value = 7558.125
$ws.Range("N45").Value = -8312.125
$ws.Range("H61").Value = 2609.8
$ws.Range("I61").Value = 2013.7894
$ws.Range("K61").Value = 2013.7894
$ws.Range("M61").Value = -1801.7894
$ws.Range("H97").Value = 624432.6
$ws.Range("I97").Value = 901466.9
$ws.Range("K97").Value = 901466.9
$ws.Range("M97").Value = -900970.9
$ws.Range("H110").Value = 3984890.5
$ws.Range("I110").Value = 4630520
$ws.Range("K110").Value = 4630520
$ws.Range("M110").Value = -4628475
$ws.Range("H119").Value = 44864
$ws.Range("J119").Value = 44864
$ws.Range("L119").Value = 44864
$ws.Range("N119").Value = -54540
$ws.Range("H122").Value = 6753901.5
$ws.Range("I122").Value = 12347105
$ws.Range("K122").Value = 37041315
$ws.Range("M122").Value = -37038865
$ws.Range("H132").Value = 2034.6
$ws.Range("I132").Value = 1846.6666
$ws.Range("J132").Value = 2786.3333
$ws.Range("K132").Value = 5539.9998
$ws.Range("L132").Value = 8358.999899999999
$ws.Range("M132").Value = -3009.9998
$ws.Range("N132").Value = -13418.9999
$ws.Range("H136").Value = 2609.8
$ws.Range("I136").Value = 2013.7894
$ws.Range("K136").Value = 6041.3682
$ws.Range("M136").Value = -3491.3682

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3845.6333
$ws.Range("I134").Value = 1451.5883
$ws.Range("K134").Value = 4354.7649
$ws.Range("M134").Value = -1819.7649

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 1872.625
$ws.Range("I122").Value = 1634
$ws.Range("J122").Value = 4497.5
$ws.Range("K122").Value = 4902
$ws.Range("L122").Value = 13492.5
$ws.Range("M122").Value = -2452
$ws.Range("N122").Value = -18392.5
$ws.Range("H132").Value = 57158.676
$ws.Range("I132").Value = 37364.895
$ws.Range("K132").Value = 112094.685
$ws.Range("M132").Value = -109564.685
$ws.Range("H134").Value = 21779.11
$ws.Range("I134").Value = 29104.031
$ws.Range("J134").Value = 3185.077
$ws.Range("K134").Value = 87312.09299999999
$ws.Range("L134").Value = 9555.231
$ws.Range("M134").Value = -84777.09299999999
$ws.Range("N134").Value = -14625.231

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 93.666664
$ws.Range("J2").Value = 128.88889
$ws.Range("L2").Value = 773.33334
$ws.Range("N2").Value = -999.33334
$ws.Range("H38").Value = 78.73333
$ws.Range("I38").Value = 38.285713
$ws.Range("J38").Value = 114.125
$ws.Range("K38").Value = 114.857139
$ws.Range("L38").Value = 342.375
$ws.Range("M38").Value = 232.142861
$ws.Range("N38").Value = -1036.375
$ws.Range("H57").Value = 3091
$ws.Range("I57").Value = 782
$ws.Range("J57").Value = 5400
$ws.Range("K57").Value = 2346
$ws.Range("L57").Value = 16200
$ws.Range("M57").Value = -1787
$ws.Range("N57").Value = -17318

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 22237590
$ws.Range("I70").Value = 25016412
$ws.Range("K70").Value = 25016412
$ws.Range("M70").Value = -25016142
$ws.Range("H73").Value = 22237590
$ws.Range("I73").Value = 25016412
$ws.Range("K73").Value = 25016412
$ws.Range("M73").Value = -25015476
$ws.Range("H113").Value = 11908632
$ws.Range("I113").Value = 41668584
$ws.Range("K113").Value = 41668584
$ws.Range("M113").Value = -41666414
$ws.Range("H121").Value = 35997.2
$ws.Range("J121").Value = 35997.2
$ws.Range("L121").Value = 35997.2
$ws.Range("N121").Value = -39491.2
$ws.Range("H122").Value = 447089.66
$ws.Range("I122").Value = 557712.0600000001
$ws.Range("J122").Value = 4600
$ws.Range("K122").Value = 1673136.18
$ws.Range("L122").Value = 13800
$ws.Range("M122").Value = -1670686.18
$ws.Range("N122").Value = -18700
$ws.Range("H126").Value = 4157976.8
$ws.Range("I126").Value = 2843886
$ws.Range("K126").Value = 8531658
$ws.Range("M126").Value = -8529188
$ws.Range("H132").Value = 2776.878
$ws.Range("J132").Value = 4194
$ws.Range("L132").Value = 12582
$ws.Range("N132").Value = -17642

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1599.2084
$ws.Range("I16").Value = 962.0526
$ws.Range("J16").Value = 4020.4
$ws.Range("K16").Value = 962.0526
$ws.Range("L16").Value = 4020.4
$ws.Range("M16").Value = -792.0526
$ws.Range("N16").Value = -4360.4
$ws.Range("H22").Value = 128083.43
$ws.Range("J22").Value = 1499.2
$ws.Range("L22").Value = 1499.2
$ws.Range("N22").Value = -2089.2
$ws.Range("H27").Value = 128083.43
$ws.Range("J27").Value = 1499.2
$ws.Range("L27").Value = 1499.2
$ws.Range("N27").Value = -1713.2
$ws.Range("H46").Value = 8370.583000000001
$ws.Range("I46").Value = 6166.6665
$ws.Range("J46").Value = 9105.223
$ws.Range("K46").Value = 6166.6665
$ws.Range("L46").Value = 9105.223
$ws.Range("M46").Value = -5978.6665
$ws.Range("N46").Value = -9481.223
$ws.Range("H55").Value = 1450.75
$ws.Range("I55").Value = 1422.4166
$ws.Range("J55").Value = 1493.25
$ws.Range("K55").Value = 1422.4166
$ws.Range("L55").Value = 1493.25
$ws.Range("M55").Value = -1249.4166
$ws.Range("N55").Value = -1839.25
$ws.Range("I61").Value = 2714980.2
$ws.Range("K61").Value = 2714980.2
$ws.Range("M61").Value = -2714778.2
$ws.Range("H68").Value = 722.25
$ws.Range("I68").Value = 722.25
$ws.Range("K68").Value = 722.25
$ws.Range("M68").Value = 26.75
$ws.Range("H71").Value = 722.25
$ws.Range("I71").Value = 722.25
$ws.Range("K71").Value = 3611.25
$ws.Range("M71").Value = 132.75
$ws.Range("I113").Value = 2714980.2
$ws.Range("K113").Value = 2714980.2
$ws.Range("M113").Value = -2712810.2
$ws.Range("H119").Value = 62500
$ws.Range("J119").Value = 62500
$ws.Range("L119").Value = 62500
$ws.Range("N119").Value = -72176
$ws.Range("H122").Value = 4725.4287
$ws.Range("I122").Value = 2896.8
$ws.Range("K122").Value = 8690.400000000001
$ws.Range("M122").Value = -6240.400000000001
$ws.Range("H136").Value = 46052.277
$ws.Range("I136").Value = 64889.438
$ws.Range("J136").Value = 5866.3335
$ws.Range("K136").Value = 194668.314
$ws.Range("L136").Value = 17599.0005
$ws.Range("M136").Value = -192118.314
$ws.Range("N136").Value = -22699.0005

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H76").Value = 0
$ws.Range("J76").Value = 0
$ws.Range("L76").Value = 0
$ws.Range("N76").ClearContents()
$ws.Range("H79").Value = 0
$ws.Range("J79").Value = 0
$ws.Range("L79").Value = 0
$ws.Range("N79").ClearContents()
$ws.Range("H113").Value = 772.96155
$ws.Range("I113").Value = 544.5454999999999
$ws.Range("J113").Value = 940.4666999999999
$ws.Range("K113").Value = 1633.6365
$ws.Range("L113").Value = 2821.4001
$ws.Range("M113").Value = 536.3635000000002
$ws.Range("N113").Value = -7161.4001
$ws.Range("H119").Value = 29500
$ws.Range("J119").Value = 29500
$ws.Range("L119").Value = 29500
$ws.Range("N119").Value = -39176
$ws.Range("H122").Value = 2418.075
$ws.Range("I122").Value = 1648.4286
$ws.Range("K122").Value = 4945.2858
$ws.Range("M122").Value = -2495.2858
$ws.Range("H124").Value = 26666
$ws.Range("J124").Value = 26666
$ws.Range("L124").Value = 26666
$ws.Range("N124").Value = -36486
$ws.Range("H132").Value = 25274336
$ws.Range("I132").Value = 31252204
$ws.Range("K132").Value = 93756612
$ws.Range("M132").Value = -93754082
